# Weekly CompStat report refresh: new crime data collected.
# Updates the report header (issue number + week-covered date range) and
# refreshes the crime-complaint statistics table (rows 15-31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header text: "Volume 31 Number 24" -> "...25", and the week-covered
#    date range 6/10/2024-6/16/2024 -> 6/17/2024-6/23/2024.
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  25"
$ws.Range("C9").Value = "Report Covering the Week  6/17/2024  Through  6/23/2024"

# ---------------------------------------------------------------------
# 2. Some cells flip between a real number and the report's placeholder
#    text ("0" / "***.*") used when a % change is undefined. Force the
#    text ones to text (so Excel doesn't silently re-parse "0" back into
#    a number) and then repaint with the donor cells' number format so
#    the visual style matches the rest of the table.
# ---------------------------------------------------------------------
$text0Cells    = @("C16", "D26", "C29", "D29", "C30", "D30")
$textStarCells = @("E26", "E29", "E30")

foreach ($c in $text0Cells) {
    $ws.Range($c).NumberFormat = "@"
    $ws.Range($c).Value = "0"
}
foreach ($c in $textStarCells) {
    $ws.Range($c).NumberFormat = "@"
    $ws.Range($c).Value = "***.*"
}

# Donor cell C15 already carries the placeholder-text style (general
# number format, right-aligned) used throughout the sheet.
$ws.Range("C15").Copy()
foreach ($c in ($text0Cells + $textStarCells)) {
    $ws.Range($c).PasteSpecial(-4122)
}

# Cells that flip the other way: placeholder text -> a real number.
$ws.Range("C23").Value = 2
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = -100
$ws.Range("C31").Value = 1

# Donor cells G15 (plain count style) / H15 (percent style) repaint the
# number format/style on the cells above after the value write.
$ws.Range("G15").Copy()
foreach ($c in @("C23", "D28", "C31")) {
    $ws.Range($c).PasteSpecial(-4122)
}
$ws.Range("H15").Copy()
$ws.Range("E28").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3. Plain numeric refreshes across rows 15-31 (counts + % changes).
# ---------------------------------------------------------------------
$ws.Range("N15").Value = -73.333333333333

$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 4
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = -73.333333333333
$ws.Range("J16").Value = 60
$ws.Range("K16").Value = 1.666666666666
$ws.Range("L16").Value = 35.555555555555
$ws.Range("M16").Value = -18.666666666666
$ws.Range("N16").Value = -83.776595744680

$ws.Range("C17").Value = 4
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = 13.333333333333
$ws.Range("I17").Value = 90
$ws.Range("J17").Value = 83
$ws.Range("K17").Value = 8.433734939759
$ws.Range("L17").Value = 42.857142857142
$ws.Range("M17").Value = 63.636363636363
$ws.Range("N17").Value = -35.714285714285

$ws.Range("D18").Value = 2
$ws.Range("F18").Value = 2
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = -83.333333333333
$ws.Range("J18").Value = 47
$ws.Range("K18").Value = -51.063829787234
$ws.Range("M18").Value = -39.473684210526
$ws.Range("N18").Value = -92.434210526315

$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 28.571428571428
$ws.Range("F19").Value = 27
$ws.Range("G19").Value = 30
$ws.Range("H19").Value = -10
$ws.Range("I19").Value = 159
$ws.Range("J19").Value = 192
$ws.Range("K19").Value = -17.1875
$ws.Range("L19").Value = -2.453987730061
$ws.Range("M19").Value = 29.268292682926
$ws.Range("N19").Value = -49.683544303797

$ws.Range("D20").Value = 8
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = -75
$ws.Range("J20").Value = 64
$ws.Range("K20").Value = -76.5625
$ws.Range("N20").Value = -92.574257425742

$ws.Range("C21").Value = 13
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = -45.833333333333
$ws.Range("F21").Value = 53
$ws.Range("G21").Value = 85
$ws.Range("H21").Value = -37.647058823529
$ws.Range("I21").Value = 352
$ws.Range("J21").Value = 450
$ws.Range("K21").Value = -21.777777777777
$ws.Range("L21").Value = -6.133333333333
$ws.Range("M21").Value = 13.915857605178
$ws.Range("N21").Value = -74.117647058823

$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 9
$ws.Range("K22").Value = -47.058823529411
$ws.Range("L22").Value = -10
$ws.Range("M22").Value = -25

$ws.Range("D23").Value = 5
$ws.Range("E23").Value = -60
$ws.Range("G23").Value = 13
$ws.Range("H23").Value = -30.769230769230
$ws.Range("I23").Value = 65
$ws.Range("J23").Value = 73
$ws.Range("K23").Value = -10.958904109589
$ws.Range("L23").Value = 25
$ws.Range("M23").Value = 51.162790697674

$ws.Range("C24").Value = 7
$ws.Range("D24").Value = 12
$ws.Range("E24").Value = -41.666666666666
$ws.Range("F24").Value = 38
$ws.Range("G24").Value = 44
$ws.Range("H24").Value = -13.636363636363
$ws.Range("I24").Value = 183
$ws.Range("J24").Value = 250
$ws.Range("K24").Value = -26.8
$ws.Range("L24").Value = -25.306122448979
$ws.Range("M24").Value = -29.069767441860

$ws.Range("D25").Value = 5
$ws.Range("E25").Value = -80
$ws.Range("G25").Value = 17
$ws.Range("H25").Value = -64.705882352941
$ws.Range("I25").Value = 30
$ws.Range("J25").Value = 89
$ws.Range("K25").Value = -66.292134831460
$ws.Range("L25").Value = -60

$ws.Range("C26").Value = 5
$ws.Range("F26").Value = 12
$ws.Range("G26").Value = 17
$ws.Range("H26").Value = -29.411764705882
$ws.Range("I26").Value = 124
$ws.Range("K26").Value = 5.982905982905
$ws.Range("L26").Value = 6.896551724137
$ws.Range("M26").Value = -15.646258503401

$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 20
$ws.Range("K28").Value = -50
$ws.Range("L28").Value = -60

$ws.Range("M29").Value = 0
$ws.Range("N29").Value = -58.333333333333

$ws.Range("M30").Value = -20
$ws.Range("N30").Value = -63.636363636363

$ws.Range("F31").Value = 1
$ws.Range("I31").Value = 7
$ws.Range("K31").Value = 600
$ws.Range("L31").Value = 0
